# Updated cryptos list on Wed Sep  6 17:45:59 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto
# table, and swaps the EnergySwap/Mantle rows (48/49) to reflect the new
# ranking order. Price values that look like plain decimal numbers are
# written with a leading apostrophe (forcing Excel to keep them as text,
# matching the source data's inline-string storage) and the cell style is
# reset to "Normal" right afterwards so no stray number-format/quote-
# prefix style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.736.25"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.627.52"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'214.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.258"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "'0.0636"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").Value = "'0.0783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "'4.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "1.625.20"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "1.852.60"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "'0.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "'62.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "25.752.70"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'4.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "'194.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "'6.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").Value = "'139.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("D28").Value = "'6.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").Value = "'15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "'1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").Value = "'3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").Value = "'3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").Value = "'1.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'0.896"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "'2.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'0.544"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("D39").Value = "1.108.77"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").Value = "'99.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").Value = "'0.799"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "1.758.85"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").Value = "0.0₆0110"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'54.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.418"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  +3.02%  "
